# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (style index 1: bold, bordered, centered)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows ---
# I column is 1 for every row except row 22 (=4).
# J column mirrors the H column value for every row except row 22 (=6).
$iValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1;
    19 = 1; 20 = 1; 21 = 1; 22 = 4; 23 = 1
}
$jValues = @{
    2 = 5; 3 = 3; 4 = 6; 5 = 7; 6 = 6; 7 = 6; 8 = 6; 9 = 4; 10 = 7;
    11 = 5; 12 = 6; 13 = 6; 14 = 5; 15 = 4; 16 = 4; 17 = 2; 18 = 5;
    19 = 5; 20 = 5; 21 = 6; 22 = 6; 23 = 2
}

for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}

Write-Output "I0/IF columns added"
